# "update alta cuenta ahorro"
# Flip the Si/No answers for the first few test cases on the Entregable1
# sheet and leave the selection parked on B8 (matches the saved sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entregable1")

$ws.Range("B2").Value = "No"   # Alta Cuenta Ahorro: Si -> No
$ws.Range("B3").Value = "No"   # Buscar Cuenta: Si -> No
$ws.Range("B4").Value = "No"   # Alta Cuenta CTS: Si -> No
$ws.Range("B6").Value = "Si"   # Alta DPF: No -> Si

[void]$ws.Range("B8").Select()
